# =====================================================================
# khl_stats_1369_ext.xlsx -- 2025-12-21 KHL results refresh
# - Append 5 new completed matches (2025-12-21) to Matches_SOG
# - Roll forward shots-on-goal aggregates in Shots_HA / Shots_Summary
#   for the 10 teams that played on 2025-12-21
# - Bump as_of_utc / build_version in Meta_ext
# =====================================================================

$wb = $excel.ActiveWorkbook

# --- Matches_SOG: append newly completed matches (rows 395-399) ---
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$wsMatches.Cells.Item(395, 1).NumberFormat = "@"   # uid stays text, like existing rows
$wsMatches.Cells.Item(395, 1).Value = "897889"
$wsMatches.Cells.Item(395, 2).Value = "2025-12-21T13:30:00"
$wsMatches.Cells.Item(395, 3).Value = "Сибирь"
$wsMatches.Cells.Item(395, 4).Value = "Амур"
$wsMatches.Cells.Item(395, 5).Value = 24
$wsMatches.Cells.Item(395, 6).Value = 34
$wsMatches.Cells.Item(395, 7).Value = "khl_text"

$wsMatches.Cells.Item(396, 1).NumberFormat = "@"   # uid stays text, like existing rows
$wsMatches.Cells.Item(396, 1).Value = "897887"
$wsMatches.Cells.Item(396, 2).Value = "2025-12-21T15:30:00"
$wsMatches.Cells.Item(396, 3).Value = "Автомобилист"
$wsMatches.Cells.Item(396, 4).Value = "Адмирал"
$wsMatches.Cells.Item(396, 5).Value = 33
$wsMatches.Cells.Item(396, 6).Value = 25
$wsMatches.Cells.Item(396, 7).Value = "khl_text"

$wsMatches.Cells.Item(397, 1).NumberFormat = "@"   # uid stays text, like existing rows
$wsMatches.Cells.Item(397, 1).Value = "897888"
$wsMatches.Cells.Item(397, 2).Value = "2025-12-21T15:00:00"
$wsMatches.Cells.Item(397, 3).Value = "Барыс"
$wsMatches.Cells.Item(397, 4).Value = "ХК Сочи"
$wsMatches.Cells.Item(397, 5).Value = 42
$wsMatches.Cells.Item(397, 6).Value = 28
$wsMatches.Cells.Item(397, 7).Value = "khl_text"

$wsMatches.Cells.Item(398, 1).NumberFormat = "@"   # uid stays text, like existing rows
$wsMatches.Cells.Item(398, 1).Value = "897891"
$wsMatches.Cells.Item(398, 2).Value = "2025-12-21T15:30:00"
$wsMatches.Cells.Item(398, 3).Value = "ЦСКА"
$wsMatches.Cells.Item(398, 4).Value = "Динамо М"
$wsMatches.Cells.Item(398, 5).Value = 27
$wsMatches.Cells.Item(398, 6).Value = 13
$wsMatches.Cells.Item(398, 7).Value = "khl_text"

$wsMatches.Cells.Item(399, 1).NumberFormat = "@"   # uid stays text, like existing rows
$wsMatches.Cells.Item(399, 1).Value = "897890"
$wsMatches.Cells.Item(399, 2).Value = "2025-12-21T17:10:00"
$wsMatches.Cells.Item(399, 3).Value = "Динамо Мн"
$wsMatches.Cells.Item(399, 4).Value = "Лада"
$wsMatches.Cells.Item(399, 5).Value = 45
$wsMatches.Cells.Item(399, 6).Value = 21
$wsMatches.Cells.Item(399, 7).Value = "khl_text"

# --- Shots_HA: refresh as_of_utc for every team, plus home/away shot
#     totals+per-game rates for teams that played on 2025-12-21 ---
$wsHA = $wb.Worksheets.Item("Shots_HA")
$wsHA.Cells.Item(2, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(3, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(3, 5).Value = 17
$wsHA.Cells.Item(3, 7).Value = 494
$wsHA.Cells.Item(3, 8).Value = 518
$wsHA.Cells.Item(3, 9).Value = 29.1
$wsHA.Cells.Item(3, 10).Value = 30.5
$wsHA.Cells.Item(4, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(4, 6).Value = 19
$wsHA.Cells.Item(4, 11).Value = 588
$wsHA.Cells.Item(4, 12).Value = 539
$wsHA.Cells.Item(4, 13).Value = 30.9
$wsHA.Cells.Item(4, 14).Value = 28.4
$wsHA.Cells.Item(5, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(6, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(6, 6).Value = 18
$wsHA.Cells.Item(6, 11).Value = 519
$wsHA.Cells.Item(6, 12).Value = 668
$wsHA.Cells.Item(6, 13).Value = 28.8
$wsHA.Cells.Item(6, 14).Value = 37.1
$wsHA.Cells.Item(7, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(7, 5).Value = 22
$wsHA.Cells.Item(7, 7).Value = 724
$wsHA.Cells.Item(7, 8).Value = 676
$wsHA.Cells.Item(7, 9).Value = 32.9
$wsHA.Cells.Item(7, 10).Value = 30.7
$wsHA.Cells.Item(8, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(8, 6).Value = 18
$wsHA.Cells.Item(8, 11).Value = 509
$wsHA.Cells.Item(8, 12).Value = 541
$wsHA.Cells.Item(8, 13).Value = 28.3
$wsHA.Cells.Item(8, 14).Value = 30.1
$wsHA.Cells.Item(9, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(9, 5).Value = 19
$wsHA.Cells.Item(9, 7).Value = 706
$wsHA.Cells.Item(9, 8).Value = 504
$wsHA.Cells.Item(9, 9).Value = 37.2
$wsHA.Cells.Item(9, 10).Value = 26.5
$wsHA.Cells.Item(10, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(11, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(11, 6).Value = 20
$wsHA.Cells.Item(11, 11).Value = 483
$wsHA.Cells.Item(11, 12).Value = 774
$wsHA.Cells.Item(11, 13).Value = 24.1
$wsHA.Cells.Item(11, 14).Value = 38.7
$wsHA.Cells.Item(12, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(13, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(14, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(15, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(16, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(17, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(18, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(18, 5).Value = 18
$wsHA.Cells.Item(18, 7).Value = 473
$wsHA.Cells.Item(18, 8).Value = 615
$wsHA.Cells.Item(18, 9).Value = 26.3
$wsHA.Cells.Item(19, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(20, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(21, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(22, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(22, 6).Value = 17
$wsHA.Cells.Item(22, 11).Value = 443
$wsHA.Cells.Item(22, 12).Value = 630
$wsHA.Cells.Item(22, 13).Value = 26.1
$wsHA.Cells.Item(22, 14).Value = 37.1
$wsHA.Cells.Item(23, 4).Value = "2025-12-21T17:10:00Z"
$wsHA.Cells.Item(23, 5).Value = 17
$wsHA.Cells.Item(23, 7).Value = 453
$wsHA.Cells.Item(23, 8).Value = 448
$wsHA.Cells.Item(23, 10).Value = 26.4

# --- Shots_Summary: refresh as_of_utc for every team, plus combined
#     shot totals+per-game rates for teams that played on 2025-12-21 ---
$wsSum = $wb.Worksheets.Item("Shots_Summary")
$wsSum.Cells.Item(2, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(3, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(3, 5).Value = 36
$wsSum.Cells.Item(3, 6).Value = 1030
$wsSum.Cells.Item(3, 7).Value = 1117
$wsSum.Cells.Item(3, 8).Value = 28.6
$wsSum.Cells.Item(3, 9).Value = 31
$wsSum.Cells.Item(4, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(4, 5).Value = 35
$wsSum.Cells.Item(4, 6).Value = 1172
$wsSum.Cells.Item(4, 7).Value = 972
$wsSum.Cells.Item(4, 8).Value = 33.5
$wsSum.Cells.Item(4, 9).Value = 27.8
$wsSum.Cells.Item(5, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(6, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(6, 5).Value = 37
$wsSum.Cells.Item(6, 6).Value = 1103
$wsSum.Cells.Item(6, 7).Value = 1316
$wsSum.Cells.Item(6, 8).Value = 29.8
$wsSum.Cells.Item(6, 9).Value = 35.6
$wsSum.Cells.Item(7, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(7, 5).Value = 38
$wsSum.Cells.Item(7, 6).Value = 1182
$wsSum.Cells.Item(7, 7).Value = 1203
$wsSum.Cells.Item(7, 8).Value = 31.1
$wsSum.Cells.Item(7, 9).Value = 31.7
$wsSum.Cells.Item(8, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(8, 5).Value = 36
$wsSum.Cells.Item(8, 6).Value = 1094
$wsSum.Cells.Item(8, 7).Value = 1018
$wsSum.Cells.Item(8, 8).Value = 30.4
$wsSum.Cells.Item(9, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(9, 5).Value = 35
$wsSum.Cells.Item(9, 6).Value = 1255
$wsSum.Cells.Item(9, 7).Value = 947
$wsSum.Cells.Item(9, 8).Value = 35.9
$wsSum.Cells.Item(9, 9).Value = 27.1
$wsSum.Cells.Item(10, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(11, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(11, 5).Value = 36
$wsSum.Cells.Item(11, 6).Value = 905
$wsSum.Cells.Item(11, 7).Value = 1348
$wsSum.Cells.Item(11, 8).Value = 25.1
$wsSum.Cells.Item(11, 9).Value = 37.4
$wsSum.Cells.Item(12, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(13, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(14, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(15, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(16, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(17, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(18, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(18, 5).Value = 37
$wsSum.Cells.Item(18, 6).Value = 973
$wsSum.Cells.Item(18, 7).Value = 1260
$wsSum.Cells.Item(18, 8).Value = 26.3
$wsSum.Cells.Item(19, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(20, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(21, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(22, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(22, 5).Value = 34
$wsSum.Cells.Item(22, 6).Value = 941
$wsSum.Cells.Item(22, 7).Value = 1167
$wsSum.Cells.Item(22, 9).Value = 34.3
$wsSum.Cells.Item(23, 4).Value = "2025-12-21T17:10:00Z"
$wsSum.Cells.Item(23, 5).Value = 37
$wsSum.Cells.Item(23, 6).Value = 967
$wsSum.Cells.Item(23, 7).Value = 981
$wsSum.Cells.Item(23, 9).Value = 26.5

# --- Meta_ext: bump as_of_utc + build_version ---
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Cells.Item(2, 2).Value = "2025-12-21T17:10:00Z"
$wsMeta.Cells.Item(2, 4).Value = 77
